$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("EnemyDB")

# Relabel the EnemyDB header row: enemyID/health/speed/itemProp -> health/speed/damage/useTime
$ws2.Range("A1").Value() = "health"
$ws2.Range("B1").Value() = "speed"
$ws2.Range("C1").Value() = "damage"
$ws2.Range("D1").Value() = "useTime"

# Make EnemyDB the active/selected sheet with a new selection
$ws2.Activate()
$ws2.Range("E6").Select()
